$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# Rename antibody names from "Acme mAb N" to "COVIC N"
$ws.Range("A2").Value = "COVIC 1"
$ws.Range("A4").Value = "COVIC 1"
$ws.Range("A5").Value = "COVIC 4"
$ws.Range("A6").Value = "COVIC 5"
$ws.Range("A7").Value = "COVIC 6"
$ws.Range("A8").Value = "COVIC 7"
$ws.Range("A9").Value = "COVIC 8"
$ws.Range("A10").Value = "COVIC 9"
$ws.Range("A11").Value = "COVIC 10"

# Fill in the previously blank qualitative measure for row 6
$ws.Range("B6").Value = "postive"

# Clear out qualitative measures for rows 8-11 (now blank in target)
$ws.Range("B8").Value = ""
$ws.Range("B9").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
